# Move "Notice u/s 94 BNSS, 2023" up onto the same line as "To," (separated
# by a center-aligned tab stop at 4680 twips / 234pt) and drop the old
# standalone heading paragraph that used to carry that text.
$d = $word.ActiveDocument

$toPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "To,") {
        $toPara = $p
        break
    }
}

if ($toPara -ne $null) {
    # Center tab stop at 4680 twips (= 234pt) so the heading lines up nicely.
    $toPara.Format.TabStops.Add(234, 1)

    $toRange = $d.Range($toPara.Range.Start, $toPara.Range.End)
    $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
           "<w:r><w:t>To,</w:t><w:tab/></w:r>" +
           "<w:r><w:rPr><w:b/></w:rPr><w:t>Notice u/s 94 BNSS, 2023</w:t></w:r>" +
           "</w:p>"
    $toRange.InsertXML($xml) | Out-Null
}

# Remove the now-redundant standalone "Notice u/s 94 BNSS, 2023" paragraph.
$noticePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Notice u/s 94 BNSS, 2023") {
        $noticePara = $p
        break
    }
}

if ($noticePara -ne $null) {
    $noticePara.Range.Delete()
}

# Every "List Paragraph" styled paragraph picks up an explicit SpaceBefore=0
# alongside its existing SpaceAfter=0.
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "List Paragraph") {
        $p.Format.SpaceBefore = 0
    }
}
